$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row text (row 1): fix typos / rename categories,
# and swap the contents of D1 and E1.
$ws.Range("A1").Value = "Years"
$ws.Range("B1").Value = "women "
$ws.Range("C1").Value = "men"
$ws.Range("D1").Value = "Transgender.people "
$ws.Range("E1").Value = "nonconforming "

# Update the active selection to D1, as recorded in the saved view state.
$ws.Range("D1").Select()
